$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11 ("Iteration 5") data update ---
# C11: Updated On date (stored as a date serial number, numFmt already applied via style)
$ws.Range("C11").Value = 43109

# D11: Planned Tasks
$ws.Range("D11").Value = 23

# E11: Actual Tasks
$ws.Range("E11").Value = 21

# F11: Task Metric = Actual / Planned
$ws.Range("F11").Formula = "=E11/D11"

# G11: Task Metrics comment (two lines of text)
$ws.Range("G11").Value = "Estimates are fairly on track and we are fairly on track.`nDue to our tight deadline, we added in additional tasks in the next iteration"

# Row 11 grows taller to fit the wrapped comment text
$ws.Rows.Item(11).RowHeight = 72.5

# --- Sheet view / selection state ---
$excel.Goto($ws.Range("A10"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("G12").Select()
